$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Insert a new row at row 16, shifting existing rows 16-21 down to 17-22.
$ws.Rows.Item(16).Insert()

# Fill in the new task row (sequence 14 - "Configurar a taxa de crescimento automática...")
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Configurar a taxa de crescimento automática dos bancos para 100MB"
$ws.Cells.Item(16, 3).Value = "Edicarlos"
$ws.Cells.Item(16, 4).Value = 5 / 1440
$ws.Cells.Item(16, 4).NumberFormat = "h:mm"

# Renumber the "Sequência" column for the rows that were pushed down.
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(20, 1).Value = 18

# Resize the Excel table (Tabela1) to include the newly inserted row.
$tbl = $ws.ListObjects.Item("Tabela1")
$tbl.Resize($ws.Range("A2:D21"))

# Update the selection to match the saved workbook state.
$ws.Application.Goto($ws.Range("A20"))

$wb.Save()
